{"js": "// Add \"Rust\" (and \"Android\") to the list of languages used in current role,\n// plus two small related wording tweaks in the same commit:\n//   1. Summary: \"enterprise mobile software testing\" -> \"enterprise mobile app testing\"\n//   2. Experience / Development Tools: \"iOS, Swift, ...\" -> \"iOS, Android, Rust, Swift, ...\"\n//   3. Experience / Development tasks: \"and macOS.\" -> \"and macOS using multiple programming languages.\"\n\nconst body = context.document.body;\n\n// 1. Summary paragraph: \"software\" -> \"app\"\nconst summaryMatches = body.search(\"enterprise mobile software testing\", { matchCase: true });\nsummaryMatches.load(\"items\");\nawait context.sync();\nif (summaryMatches.items.length > 0) {\n  summaryMatches.items[0].insertText(\"enterprise mobile app testing\", \"Replace\");\n}\n\n// 2. Development Tools paragraph: insert \"Android, Rust, \" after \"iOS, \"\nconst toolsMatches = body.search(\"iOS, Swift, Objective-C, Xcode\", { matchCase: true });\ntoolsMatches.load(\"items\");\nawait context.sync();\nif (toolsMatches.items.length > 0) {\n  toolsMatches.items[0].insertText(\"iOS, Android, Rust, Swift, Objective-C, Xcode\", \"Replace\");\n}\n\n// 3. Development tasks paragraph: \"and macOS.\" -> \"and macOS using multiple programming languages.\"\nconst macosMatches = body.search(\"and macOS. Integrated components\", { matchCase: true });\nmacosMatches.load(\"items\");\nawait context.sync();\nif (macosMatches.items.length > 0) {\n  macosMatches.items[0].insertText(\n    \"and macOS using multiple programming languages. Integrated components\",\n    \"Replace\"\n  );\n}\n\nawait context.sync();\n", "ps1": "# Add \"Rust\" (and \"Android\") to the list of languages used in current role,\n# plus two small related wording tweaks in the same commit:\n#   1. Summary: \"enterprise mobile software testing\" -> \"enterprise mobile app testing\"\n#   2. Experience / Development Tools: \"iOS, Swift, ...\" -> \"iOS, Android, Rust, Swift, ...\"\n#   3. Experience / Development tasks: \"and macOS.\" -> \"and macOS using multiple programming languages.\"\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue=1, wdReplaceOne=1 -> replace just the (unique) first match\n    $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 1)\n}\n\n# 1. Summary paragraph: \"software\" -> \"app\"\nReplace-Text \"enterprise mobile software testing\" \"enterprise mobile app testing\"\n\n# 2. Development Tools paragraph: insert \"Android, Rust, \" after \"iOS, \"\nReplace-Text \"iOS, Swift, Objective-C, Xcode\" \"iOS, Android, Rust, Swift, Objective-C, Xcode\"\n\n# 3. Development tasks paragraph: \"and macOS.\" -> \"and macOS using multiple programming languages.\"\nReplace-Text \"and macOS. Integrated components\" \"and macOS using multiple programming languages. Integrated components\"\n"}
